$wb = $excel.ActiveWorkbook

# --- Campaign sheet: add CampaignName/LOB/IndustryGroups/HLSubGroup columns ---
$wsCampaign = $wb.Worksheets.Item("Campaign")

$wsCampaign.Range("B1").Value = "CampaignName"
$wsCampaign.Range("B2").Value = "TestParentCampaign"
$wsCampaign.Range("C1").Value = "LOB"
$wsCampaign.Range("C2").Value = "CF"
$wsCampaign.Range("D1").Value = "IndustryGroups"
$wsCampaign.Range("D2").Value = "BUS"
$wsCampaign.Range("E1").Value = "HLSubGroup"
$wsCampaign.Range("E2").Value = "CM"

# Bold the new header cells, matching the style of the existing header row
$wsCampaign.Range("B1:E1").Font.Bold = $true

# Widen the columns whose content doesn't fit the default width (AutoFit-like)
$wsCampaign.Range("B:B").ColumnWidth = 17.1
$wsCampaign.Range("D:D").ColumnWidth = 13.25

$wsCampaign.Activate() | Out-Null
$wsCampaign.Range("E2").Select() | Out-Null

# --- Activity sheet: add Campaigns column ---
$wsActivity = $wb.Worksheets.Item("Activity")

$wsActivity.Range("I2").Value = "TestParentCampaign"
$wsActivity.Range("I1").Value = "Campaigns"
$wsActivity.Range("I1").Font.Bold = $true
$wsActivity.Range("I1").HorizontalAlignment = -4108

$wsActivity.Range("I:I").ColumnWidth = 17.1

$wsActivity.Activate() | Out-Null
$wsActivity.Range("I8").Select() | Out-Null
